$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 46, shifting existing rows 46-61 down to 49-64
$ws.Rows("46:48").Insert()

# The insert copies stray formatting into A48/B48/C48/D48/H48 from the row
# above; clear it so row 48 only keeps its original E48/I48 formatting.
$ws.Range("A48:D48").Clear()
$ws.Range("H48").Clear()

# Fill in the two new data rows (46 and 47) for mouse 1350, "bunnytop high res"
$ws.Cells.Item(46, 1).Value = 1350
$ws.Cells.Item(46, 2).Value = 220225
$ws.Cells.Item(46, 3).Value = "V1"
$ws.Cells.Item(46, 4).Value = 200
$ws.Cells.Item(46, 5).Value = "002"
$ws.Cells.Item(46, 8).Value = "bunnytop high res"
$ws.Cells.Item(46, 9).Value = "6s"

$ws.Cells.Item(47, 1).Value = 1350
$ws.Cells.Item(47, 2).Value = 220225
$ws.Cells.Item(47, 3).Value = "V1"
$ws.Cells.Item(47, 4).Value = 200
$ws.Cells.Item(47, 5).Value = "003"
$ws.Cells.Item(47, 8).Value = "bunnytop high res"
$ws.Cells.Item(47, 9).Value = "6s"

$ws.Range("K47").Select()
